$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 33) mirroring the existing table structure
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Update the visible selection/scroll position as seen after the edit
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E29").Select()
